$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.305367469787598
$ws.Range("B1").Value = 3.767757654190063
$ws.Range("C1").Value = 3.983294010162354
$ws.Range("D1").Value = 2.908904790878296
$ws.Range("E1").Value = 1.049429059028625
